$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion rate text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = $ws1.Range("A1").Value()
$text = $text.Replace("1000 Bs = 4.08 = 15672.77 pesos", "1000 Bs = 3.92 = 15073.17 pesos")
$text = $text.Replace("15672.77 pesos = 4.06 = 972.76 Bs", "15073.17 pesos = 3.9 = 927.16 Bs")
$ws1.Range("A1").Value = $text

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 254.89
$ws2.Range("O10").Value = 3842
$ws2.Range("N12").Value = 3869.27
$ws2.Range("O12").Value = 238
